# Update cryptos list with latest prices/volumes scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.414.58"
$ws.Range("E2").Value = "  +2.69%  "
# Row 3
$ws.Range("D3").Value = "2.513.57"
$ws.Range("E3").Value = "  +1.06%  "
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").Value = "'593.09"
$ws.Range("E5").Value = "  +1.23%  "
# Row 6
$ws.Range("D6").Value = "'176.21"
$ws.Range("E6").Value = "  +0.02%  "
# Row 7
$ws.Range("E7").Value = "  -0.05%  "
# Row 8
$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  +0.35%  "
# Row 9
$ws.Range("D9").Value = "2.512.98"
$ws.Range("E9").Value = "  +1.06%  "
# Row 10
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  +9.57%  "
# Row 11
$ws.Range("E11").Value = "  -0.47%  "
# Row 12
$ws.Range("D12").Value = "'4.99"
$ws.Range("E12").Value = "  +1.32%  "
# Row 13
$ws.Range("D13").Value = "'0.338"
$ws.Range("E13").Value = "  -0.10%  "
# Row 15
$ws.Range("D15").Value = "'25.79"
$ws.Range("E15").Value = "  +0.61%  "
# Row 16
$ws.Range("D16").Value = "69.209.40"
$ws.Range("E16").Value = "  +2.47%  "
# Row 17
$ws.Range("D17").Value = "'0.0000175"
$ws.Range("E17").Value = "  +2.32%  "
# Row 18
$ws.Range("D18").Value = "2.508.31"
$ws.Range("E18").Value = "  +0.06%  "
# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  +0.53%  "
# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'361.86"
$ws.Range("E20").Value = "  +3.31%  "
# Row 21
$ws.Range("D21").Value = "'10.96"
$ws.Range("E21").Value = "  +0.49%  "
# Row 22
$ws.Range("D22").Value = "'4.03"
$ws.Range("E22").Value = "  -1.50%  "
# Row 23
$ws.Range("E23").Value = "  +0.01%  "
# Row 24
$ws.Range("D24").Value = "'70.18"
$ws.Range("E24").Value = "  -0.60%  "
# Row 25
$ws.Range("D25").Value = "'4.19"
$ws.Range("E25").Value = "  -1.53%  "
# Row 26
$ws.Range("D26").Value = "'8.99"
$ws.Range("E26").Value = "  -1.29%  "
# Row 27
$ws.Range("D27").Value = "'1.66"
$ws.Range("E27").Value = "  -5.08%  "
# Row 28
$ws.Range("D28").Value = "2.637.19"
$ws.Range("E28").Value = "  +0.84%  "
# Row 29
$ws.Range("E29").Value = "  +0.50%  "
# Row 30
$ws.Range("D30").Value = "'507.89"
$ws.Range("E30").Value = "  +0.26%  "
# Row 31
$ws.Range("D31").Value = "0.0₃0884"
$ws.Range("E31").Value = "  -1.77%  "
# Row 32
$ws.Range("D32").Value = "'7.72"
$ws.Range("E32").Value = "  -0.92%  "
# Row 33
$ws.Range("D33").Value = "'1.22"
$ws.Range("E33").Value = "  -2.84%  "
# Row 34
$ws.Range("D34").Value = "'1.77"
$ws.Range("E34").Value = "  +0.26%  "
# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.07%  "
# Row 36
$ws.Range("D36").Value = "'161.75"
$ws.Range("E36").Value = "  -0.33%  "
# Row 37
$ws.Range("E37").Value = "  -2.69%  "
# Row 38
$ws.Range("D38").Value = "'18.65"
$ws.Range("E38").Value = "  +1.96%  "
# Row 39
$ws.Range("D39").Value = "'18.70"
$ws.Range("E39").Value = "  +0.14%  "
# Row 40
$ws.Range("E40").Value = "  +0.09%  "
# Row 41
$ws.Range("D41").Value = "'1.30"
$ws.Range("E41").Value = "  -2.50%  "
# Row 42
$ws.Range("D42").Value = "'1.70"
$ws.Range("E42").Value = "  -2.28%  "
# Row 43
$ws.Range("D43").Value = "'4.77"
$ws.Range("E43").Value = "  -1.19%  "
# Row 44
$ws.Range("D44").Value = "'0.319"
$ws.Range("E44").Value = "  -2.60%  "
# Row 45
$ws.Range("D45").Value = "'2.30"
$ws.Range("E45").Value = "  -4.50%  "
# Row 46
$ws.Range("D46").Value = "'149.57"
$ws.Range("E46").Value = "  +3.37%  "
# Row 47
$ws.Range("D47").Value = "'3.55"
$ws.Range("E47").Value = "  +1.04%  "
# Row 48
$ws.Range("D48").Value = "'0.512"
$ws.Range("E48").Value = "  -0.30%  "
# Row 49
$ws.Range("D49").Value = "'0.0735"
$ws.Range("E49").Value = "  -0.89%  "
# Row 50
$ws.Range("E50").Value = "  -1.76%  "
# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0245"
$ws.Range("E51").Value = "  -3.91%  "
